$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("G2").Value = 2.3842205
$ws.Range("H2").Value = 4.768441
$ws.Range("I2").Value = 0.0684902599354226
$ws.Range("J2").Value = 0.05735520746201143
$ws.Range("K2").Value = 2
$ws.Range("M2").Value = 44.544241
$ws.Range("N2").Value = 89.088482
$ws.Range("O2").Value = 0.2677188803968527
$ws.Range("P2").Value = 0.2025351964154738
$ws.Range("Q2").Value = 106.2032925491405
$ws.Range("R2").Value = 424.813170196562
$ws.Range("S2").Value = 0.01833613570800075
$ws.Range("T2").Value = 0.01161644820876873

$ws.Range("E3").Value = 2
$ws.Range("G3").Value = 2.3842205
$ws.Range("H3").Value = 4.768441
$ws.Range("I3").Value = 0.0684902599354226
$ws.Range("J3").Value = 0.05735520746201143
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 20.36664933333333
$ws.Range("N3").Value = 61.099948
$ws.Range("O3").Value = 0.1224072166131488
$ws.Range("P3").Value = 0.1389056103700951
$ws.Range("Q3").Value = 48.55858285684467
$ws.Range("R3").Value = 291.351497141068
$ws.Range("S3").Value = 0.008383702083806142
$ws.Range("T3").Value = 0.007966960100414132

$ws.Range("E4").Value = 2
$ws.Range("G4").Value = 2.3842205
$ws.Range("H4").Value = 4.768441
$ws.Range("I4").Value = 0.0684902599354226
$ws.Range("J4").Value = 0.05735520746201143
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 25.30234
$ws.Range("N4").Value = 75.90701999999999
$ws.Range("O4").Value = 0.1520716030658262
$ws.Range("P4").Value = 0.1725682474308328
$ws.Range("Q4").Value = 60.32635772597
$ws.Range("R4").Value = 361.95814635582
$ws.Range("S4").Value = 0.01041542362277484
$ws.Range("T4").Value = 0.009897687632751137

$ws.Range("E5").Value = 2
$ws.Range("G5").Value = 2.3842205
$ws.Range("H5").Value = 4.768441
$ws.Range("I5").Value = 0.0684902599354226
$ws.Range("J5").Value = 0.05735520746201143
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 36.81412633333333
$ws.Range("N5").Value = 110.442379
$ws.Range("O5").Value = 0.2212595043374584
$ws.Range("P5").Value = 0.2510814913577403
$ws.Range("Q5").Value = 87.77299469352317
$ws.Range("R5").Value = 526.637968161139
$ws.Range("S5").Value = 0.01515412096525529
$ws.Range("T5").Value = 0.01440083102669442

$ws.Range("E6").Value = 2
$ws.Range("G6").Value = 2.3842205
$ws.Range("H6").Value = 4.768441
$ws.Range("I6").Value = 0.0684902599354226
$ws.Range("J6").Value = 0.05735520746201143
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 24.61478433333333
$ws.Range("N6").Value = 73.844353
$ws.Range("O6").Value = 0.1479392701500961
$ws.Range("P6").Value = 0.1678789468994273
$ws.Range("Q6").Value = 58.68707341061216
$ws.Range("R6").Value = 352.122440463673
$ws.Range("S6").Value = 0.01013239906723679
$ws.Range("T6").Value = 0.00962873182792065

$ws.Range("E7").Value = 2
$ws.Range("G7").Value = 2.3842205
$ws.Range("H7").Value = 4.768441
$ws.Range("I7").Value = 0.0684902599354226
$ws.Range("J7").Value = 0.05735520746201143
$ws.Range("K7").Value = 2
$ws.Range("M7").Value = 14.742243
$ws.Range("N7").Value = 29.484486
$ws.Range("O7").Value = 0.0886035254366179
$ws.Range("P7").Value = 0.06703050752643071
$ws.Range("Q7").Value = 35.1487579765815
$ws.Range("R7").Value = 140.595031906326
$ws.Range("S7").Value = 0.006068478488348788
$ws.Range("T7").Value = 0.003844548665462352

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 11.717184
$ws.Range("H8").Value = 35.151552
$ws.Range("I8").Value = 0.3365934391853332
$ws.Range("J8").Value = 0.4228058096077277
$ws.Range("K8").Value = 2
$ws.Range("M8").Value = 44.544241
$ws.Range("N8").Value = 89.088482
$ws.Range("O8").Value = 0.2677188803968527
$ws.Range("P8").Value = 0.2025351964154738
$ws.Range("Q8").Value = 521.933067937344
$ws.Range("R8").Value = 3131.598407624064
$ws.Range("S8").Value = 0.09011241868762353
$ws.Range("T8").Value = 0.08563305769450455

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 11.717184
$ws.Range("H9").Value = 35.151552
$ws.Range("I9").Value = 0.3365934391853332
$ws.Range("J9").Value = 0.4228058096077277
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 20.36664933333333
$ws.Range("N9").Value = 61.099948
$ws.Range("O9").Value = 0.1224072166131488
$ws.Range("P9").Value = 0.1389056103700951
$ws.Range("Q9").Value = 238.639777702144
$ws.Range("R9").Value = 2147.757999319296
$ws.Range("S9").Value = 0.04120146602092382
$ws.Range("T9").Value = 0.05873009905158365

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 11.717184
$ws.Range("H10").Value = 35.151552
$ws.Range("I10").Value = 0.3365934391853332
$ws.Range("J10").Value = 0.4228058096077277
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 25.30234
$ws.Range("N10").Value = 75.90701999999999
$ws.Range("O10").Value = 0.1520716030658262
$ws.Range("P10").Value = 0.1725682474308328
$ws.Range("Q10").Value = 296.47217341056
$ws.Range("R10").Value = 2668.24956069504
$ws.Range("S10").Value = 0.0511863038783533
$ws.Range("T10").Value = 0.07296285756757995

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 11.717184
$ws.Range("H11").Value = 35.151552
$ws.Range("I11").Value = 0.3365934391853332
$ws.Range("J11").Value = 0.4228058096077277
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 36.81412633333333
$ws.Range("N11").Value = 110.442379
$ws.Range("O11").Value = 0.2212595043374584
$ws.Range("P11").Value = 0.2510814913577403
$ws.Range("Q11").Value = 431.3578920469121
$ws.Range("R11").Value = 3882.221028422208
$ws.Range("S11").Value = 0.07447449751738727
$ws.Range("T11").Value = 0.1061587132310251

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 11.717184
$ws.Range("H12").Value = 35.151552
$ws.Range("I12").Value = 0.3365934391853332
$ws.Range("J12").Value = 0.4228058096077277
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 24.61478433333333
$ws.Range("N12").Value = 73.844353
$ws.Range("O12").Value = 0.1479392701500961
$ws.Range("P12").Value = 0.1678789468994273
$ws.Range("Q12").Value = 288.415957153984
$ws.Range("R12").Value = 2595.743614385856
$ws.Range("S12").Value = 0.04979538773038897
$ws.Range("T12").Value = 0.07098019405990508

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 11.717184
$ws.Range("H13").Value = 35.151552
$ws.Range("I13").Value = 0.3365934391853332
$ws.Range("J13").Value = 0.4228058096077277
$ws.Range("K13").Value = 2
$ws.Range("M13").Value = 14.742243
$ws.Range("N13").Value = 29.484486
$ws.Range("O13").Value = 0.0886035254366179
$ws.Range("P13").Value = 0.06703050752643071
$ws.Range("Q13").Value = 172.737573803712
$ws.Range("R13").Value = 1036.425442822272
$ws.Range("S13").Value = 0.02982336535065637
$ws.Range("T13").Value = 0.02834088800312942

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 1.799402
$ws.Range("H14").Value = 5.398206
$ws.Range("I14").Value = 0.05169048362276865
$ws.Range("J14").Value = 0.06493007359274758
$ws.Range("K14").Value = 2
$ws.Range("M14").Value = 44.544241
$ws.Range("N14").Value = 89.088482
$ws.Range("O14").Value = 0.2677188803968527
$ws.Range("P14").Value = 0.2025351964154738
$ws.Range("Q14").Value = 80.15299634388199
$ws.Range("R14").Value = 480.917978063292
$ws.Range("S14").Value = 0.01383851840265947
$ws.Range("T14").Value = 0.0131506252083783

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 1.799402
$ws.Range("H15").Value = 5.398206
$ws.Range("I15").Value = 0.05169048362276865
$ws.Range("J15").Value = 0.06493007359274758
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 20.36664933333333
$ws.Range("N15").Value = 61.099948
$ws.Range("O15").Value = 0.1224072166131488
$ws.Range("P15").Value = 0.1389056103700951
$ws.Range("Q15").Value = 36.64778954369866
$ws.Range("R15").Value = 329.830105893288
$ws.Range("S15").Value = 0.006327288225650664
$ws.Range("T15").Value = 0.009019151503775798

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 1.799402
$ws.Range("H16").Value = 5.398206
$ws.Range("I16").Value = 0.05169048362276865
$ws.Range("J16").Value = 0.06493007359274758
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 25.30234
$ws.Range("N16").Value = 75.90701999999999
$ws.Range("O16").Value = 0.1520716030658262
$ws.Range("P16").Value = 0.1725682474308328
$ws.Range("Q16").Value = 45.52908120068
$ws.Range("R16").Value = 409.7617308061199
$ws.Range("S16").Value = 0.007860654707762264
$ws.Range("T16").Value = 0.01120486900545545

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 1.799402
$ws.Range("H17").Value = 5.398206
$ws.Range("I17").Value = 0.05169048362276865
$ws.Range("J17").Value = 0.06493007359274758
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 36.81412633333333
$ws.Range("N17").Value = 110.442379
$ws.Range("O17").Value = 0.2212595043374584
$ws.Range("P17").Value = 0.2510814913577403
$ws.Range("Q17").Value = 66.24341255245267
$ws.Range("R17").Value = 596.190712972074
$ws.Range("S17").Value = 0.0114370107853373
$ws.Range("T17").Value = 0.01630273971163489

$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 1.799402
$ws.Range("H18").Value = 5.398206
$ws.Range("I18").Value = 0.05169048362276865
$ws.Range("J18").Value = 0.06493007359274758
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 24.61478433333333
$ws.Range("N18").Value = 73.844353
$ws.Range("O18").Value = 0.1479392701500961
$ws.Range("P18").Value = 0.1678789468994273
$ws.Range("Q18").Value = 44.29189215896866
$ws.Range("R18").Value = 398.627029430718
$ws.Range("S18").Value = 0.007647052420857892
$ws.Range("T18").Value = 0.01090039237685277

$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 1.799402
$ws.Range("H19").Value = 5.398206
$ws.Range("I19").Value = 0.05169048362276865
$ws.Range("J19").Value = 0.06493007359274758
$ws.Range("K19").Value = 2
$ws.Range("M19").Value = 14.742243
$ws.Range("N19").Value = 29.484486
$ws.Range("O19").Value = 0.0886035254366179
$ws.Range("P19").Value = 0.06703050752643071
$ws.Range("Q19").Value = 26.527221538686
$ws.Range("R19").Value = 159.163329232116
$ws.Range("S19").Value = 0.004579959080501063
$ws.Range("T19").Value = 0.004352295786650366

$ws.Range("E20").Value = 2
$ws.Range("G20").Value = 18.910282
$ws.Range("H20").Value = 37.820564
$ws.Range("I20").Value = 0.5432258172564757
$ws.Range("J20").Value = 0.4549089093375132
$ws.Range("K20").Value = 2
$ws.Range("M20").Value = 44.544241
$ws.Range("N20").Value = 89.088482
$ws.Range("O20").Value = 0.2677188803968527
$ws.Range("P20").Value = 0.2025351964154738
$ws.Range("Q20").Value = 842.3441587859619
$ws.Range("R20").Value = 3369.376635143848
$ws.Range("S20").Value = 0.1454318075985689
$ws.Range("T20").Value = 0.09213506530382219

$ws.Range("E21").Value = 2
$ws.Range("G21").Value = 18.910282
$ws.Range("H21").Value = 37.820564
$ws.Range("I21").Value = 0.5432258172564757
$ws.Range("J21").Value = 0.4549089093375132
$ws.Range("K21").Value = 3
$ws.Range("M21").Value = 20.36664933333333
$ws.Range("N21").Value = 61.099948
$ws.Range("O21").Value = 0.1224072166131488
$ws.Range("P21").Value = 0.1389056103700951
$ws.Range("Q21").Value = 385.1390822884453
$ws.Range("R21").Value = 2310.834493730672
$ws.Range("S21").Value = 0.06649476028276821
$ws.Range("T21").Value = 0.06318939971432154

$ws.Range("E22").Value = 2
$ws.Range("G22").Value = 18.910282
$ws.Range("H22").Value = 37.820564
$ws.Range("I22").Value = 0.5432258172564757
$ws.Range("J22").Value = 0.4549089093375132
$ws.Range("K22").Value = 3
$ws.Range("M22").Value = 25.30234
$ws.Range("N22").Value = 75.90701999999999
$ws.Range("O22").Value = 0.1520716030658262
$ws.Range("P22").Value = 0.1725682474308328
$ws.Range("Q22").Value = 478.4743846598799
$ws.Range("R22").Value = 2870.846307959279
$ws.Range("S22").Value = 0.0826092208569358
$ws.Range("T22").Value = 0.07850283322504627

$ws.Range("E23").Value = 2
$ws.Range("G23").Value = 18.910282
$ws.Range("H23").Value = 37.820564
$ws.Range("I23").Value = 0.5432258172564757
$ws.Range("J23").Value = 0.4549089093375132
$ws.Range("K23").Value = 3
$ws.Range("M23").Value = 36.81412633333333
$ws.Range("N23").Value = 110.442379
$ws.Range("O23").Value = 0.2212595043374584
$ws.Range("P23").Value = 0.2510814913577403
$ws.Range("Q23").Value = 696.1655105469594
$ws.Range("R23").Value = 4176.993063281756
$ws.Range("S23").Value = 0.1201938750694785
$ws.Range("T23").Value = 0.1142192073883859

$ws.Range("E24").Value = 2
$ws.Range("G24").Value = 18.910282
$ws.Range("H24").Value = 37.820564
$ws.Range("I24").Value = 0.5432258172564757
$ws.Range("J24").Value = 0.4549089093375132
$ws.Range("K24").Value = 3
$ws.Range("M24").Value = 24.61478433333333
$ws.Range("N24").Value = 73.844353
$ws.Range("O24").Value = 0.1479392701500961
$ws.Range("P24").Value = 0.1678789468994273
$ws.Range("Q24").Value = 465.4725131125153
$ws.Range("R24").Value = 2792.835078675092
$ws.Range("S24").Value = 0.08036443093161251
$ws.Range("T24").Value = 0.07636962863474873

$ws.Range("E25").Value = 2
$ws.Range("G25").Value = 18.910282
$ws.Range("H25").Value = 37.820564
$ws.Range("I25").Value = 0.5432258172564757
$ws.Range("J25").Value = 0.4549089093375132
$ws.Range("K25").Value = 2
$ws.Range("M25").Value = 14.742243
$ws.Range("N25").Value = 29.484486
$ws.Range("O25").Value = 0.0886035254366179
$ws.Range("P25").Value = 0.06703050752643071
$ws.Range("Q25").Value = 278.7799724425259
$ws.Range("R25").Value = 1115.119889770104
$ws.Range("S25").Value = 0.04813172251711169
$ws.Range("T25").Value = 0.03049277507118856
